# Commit: "messed up some spacing"
#
# Fixes the position/size of three shapes on slide 1 of the poster:
#   - Shape id=10 "Content Placeholder 2" (Purpose box)      -> widen
#   - Shape id=13 "Group 12"                                  -> reposition
#   - Shape id=28 "Content Placeholder 2" (caption box)       -> reposition + narrow
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points, while the
# OOXML stores EMU (1 pt = 12700 EMU). To land on an exact EMU value we
# compute, for each target EMU, a point value whose internal float32
# round-trip reproduces that EMU exactly.

function Find-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape id=10: "Content Placeholder 2" -- widen (ext cx 3842333 -> 4045362 EMU) ---
$shape10 = Find-ShapeById $s 10
$shape10.Width = 318.532470703125

# --- Shape id=13: "Group 12" -- move
#     (off x 8085086 -> 8222896 EMU, off y 4221575 -> 4222875 EMU) ---
$group13 = Find-ShapeById $s 13
$group13.Left = 647.47216796875
$group13.Top  = 332.5098571777344

# --- Shape id=28: "Content Placeholder 2" -- move + narrow
#     (off x 9037317 -> 9055224 EMU, off y 6583214 -> 6581319 EMU,
#      ext cx 2238685 -> 2228189 EMU) ---
$shape28 = Find-ShapeById $s 28
$shape28.Left  = 713.009765625
$shape28.Top   = 518.214111328125
$shape28.Width = 175.44796752929688
